$p = $ppt.ActivePresentation

# Slide 11 (ArrayIndex) - three TextBox shapes get resized/repositioned and
# re-fonted to 11pt Consolas.
$s = $p.Slides.Item(11)

# EMU -> point helper (1 pt = 12700 EMU), matching the OOXML units in the diff.
function EmuToPt([double]$emu) { return $emu / 12700.0 }

# --- TextBox 9 ("Possible max and min length") ---
$tb9 = $s.Shapes.Item(7)
$tb9.Left   = EmuToPt 2865898
$tb9.Top    = EmuToPt 2587450
$tb9.Width  = EmuToPt 2262158
$tb9.Height = EmuToPt 261610
$tr9 = $tb9.TextFrame.TextRange
$tr9.Font.Size = 11
$tr9.Font.NameAscii = "Consolas"
$tr9.Font.NameComplexScript = "Consolas"

# --- TextBox 11 ("Possible prefixes") ---
$tb11 = $s.Shapes.Item(9)
$tb11.Left   = EmuToPt 4837562
$tb11.Top    = EmuToPt 1638528
$tb11.Width  = EmuToPt 1492716
$tb11.Height = EmuToPt 261610
$tr11 = $tb11.TextFrame.TextRange
$tr11.Font.Size = 11
$tr11.Font.NameAscii = "Consolas"
$tr11.Font.NameComplexScript = "Consolas"

# --- TextBox 13 ("Possible string contains evaluated dynamically") ---
$tb13 = $s.Shapes.Item(11)
$tb13.Left   = EmuToPt 5289565
$tb13.Top    = EmuToPt 2605957
$tb13.Width  = EmuToPt 3689336
$tb13.Height = EmuToPt 261610
$tr13 = $tb13.TextFrame.TextRange
$tr13.Font.Size = 11
$tr13.Font.NameAscii = "Consolas"
$tr13.Font.NameComplexScript = "Consolas"
